# PerformanceTracker.xlsx - "More work from 12/6/16"
# Add round-3 results (lm, gbm, rf, MARS) below the existing round-1/round-2 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: lm - round 3, only Train RMSE / Test RMSE recorded so far
$ws.Range("A10").Value = "lm"
$ws.Range("B10").Value = 3
$ws.Range("D10").Value = 0.1019815
$ws.Range("E10").Value = 0.1297403

# Row 11: gbm - round 3, just queued (no metrics yet)
$ws.Range("A11").Value = "gbm"
$ws.Range("B11").Value = 3

# Row 12: rf - round 3, just queued (no metrics yet)
$ws.Range("A12").Value = "rf"
$ws.Range("B12").Value = 3

# Row 13: MARS - round 3, just queued (no metrics yet)
$ws.Range("A13").Value = "MARS"
$ws.Range("B13").Value = 3

# Scroll the window down a bit and leave the selection on E11, matching
# where the author was working when they saved.
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E11").Select()
